$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(2, 2).Value = "Ccl12"
$ws.Cells.Item(2, 3).Value = "Ccr10"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 42.70362466666668
$ws.Cells.Item(2, 8).Value = 128.110874
$ws.Cells.Item(2, 9).Value = 0.510021191154308
$ws.Cells.Item(2, 10).Value = 0.5102913077099245
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 0.2651005
$ws.Cells.Item(2, 14).Value = 0.530201
$ws.Cells.Item(2, 15).Value = 0.2682690254597046
$ws.Cells.Item(2, 16).Value = 0.2682690254597046
$ws.Cells.Item(2, 17).Value = 11.32075225094567
$ws.Cells.Item(2, 18).Value = 67.92451350567401
$ws.Cells.Item(2, 19).Value = 0.1368228879147639
$ws.Cells.Item(2, 20).Value = 0.1368953518198997

$ws.Cells.Item(3, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 2).Value = "Ccl12"
$ws.Cells.Item(3, 3).Value = "Ccr10"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 42.70362466666668
$ws.Cells.Item(3, 8).Value = 128.110874
$ws.Cells.Item(3, 9).Value = 0.510021191154308
$ws.Cells.Item(3, 10).Value = 0.5102913077099245
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.5
$ws.Cells.Item(3, 13).Value = 0.7230885
$ws.Cells.Item(3, 14).Value = 1.446177
$ws.Cells.Item(3, 15).Value = 0.7317309745402955
$ws.Cells.Item(3, 16).Value = 0.7317309745402955
$ws.Cells.Item(3, 17).Value = 30.87849990478301
$ws.Cells.Item(3, 18).Value = 185.270999428698
$ws.Cells.Item(3, 19).Value = 0.3731983032395441
$ws.Cells.Item(3, 20).Value = 0.3733959558900248

$ws.Cells.Item(4, 1).Value = "MuSCs"
$ws.Cells.Item(4, 2).Value = "Ccl12"
$ws.Cells.Item(4, 3).Value = "Ccr10"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.5
$ws.Cells.Item(4, 7).Value = 0.132963
$ws.Cells.Item(4, 8).Value = 0.265926
$ws.Cells.Item(4, 9).Value = 0.001588013855235666
$ws.Cells.Item(4, 10).Value = 0.001059236597621443
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.5
$ws.Cells.Item(4, 13).Value = 0.2651005
$ws.Cells.Item(4, 14).Value = 0.530201
$ws.Cells.Item(4, 15).Value = 0.2682690254597046
$ws.Cells.Item(4, 16).Value = 0.2682690254597046
$ws.Cells.Item(4, 17).Value = 0.0352485577815
$ws.Cells.Item(4, 18).Value = 0.140994231126
$ws.Cells.Item(4, 19).Value = 0.0004260149293605806
$ws.Cells.Item(4, 20).Value = 0.0002841603697751577

$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Ccl12"
$ws.Cells.Item(5, 3).Value = "Ccr10"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.5
$ws.Cells.Item(5, 7).Value = 0.132963
$ws.Cells.Item(5, 8).Value = 0.265926
$ws.Cells.Item(5, 9).Value = 0.001588013855235666
$ws.Cells.Item(5, 10).Value = 0.001059236597621443
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.5
$ws.Cells.Item(5, 13).Value = 0.7230885
$ws.Cells.Item(5, 14).Value = 1.446177
$ws.Cells.Item(5, 15).Value = 0.7317309745402955
$ws.Cells.Item(5, 16).Value = 0.7317309745402955
$ws.Cells.Item(5, 17).Value = 0.0961440162255
$ws.Cells.Item(5, 18).Value = 0.384576064902
$ws.Cells.Item(5, 19).Value = 0.001161998925875086
$ws.Cells.Item(5, 20).Value = 0.0007750762278462852

$ws.Cells.Item(6, 1).Value = "Neutrophils"
$ws.Cells.Item(6, 2).Value = "Ccl12"
$ws.Cells.Item(6, 3).Value = "Ccr10"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 8.925702
$ws.Cells.Item(6, 8).Value = 26.777106
$ws.Cells.Item(6, 9).Value = 0.1066021257320059
$ws.Cells.Item(6, 10).Value = 0.1066585841684857
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.5
$ws.Cells.Item(6, 13).Value = 0.2651005
$ws.Cells.Item(6, 14).Value = 0.530201
$ws.Cells.Item(6, 15).Value = 0.2682690254597046
$ws.Cells.Item(6, 16).Value = 0.2682690254597046
$ws.Cells.Item(6, 17).Value = 2.366208063051
$ws.Cells.Item(6, 18).Value = 14.197248378306
$ws.Cells.Item(6, 19).Value = 0.02859804838205812
$ws.Cells.Item(6, 20).Value = 0.02861319443179153

$ws.Cells.Item(7, 1).Value = "Neutrophils"
$ws.Cells.Item(7, 2).Value = "Ccl12"
$ws.Cells.Item(7, 3).Value = "Ccr10"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 8.925702
$ws.Cells.Item(7, 8).Value = 26.777106
$ws.Cells.Item(7, 9).Value = 0.1066021257320059
$ws.Cells.Item(7, 10).Value = 0.1066585841684857
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.5
$ws.Cells.Item(7, 13).Value = 0.7230885
$ws.Cells.Item(7, 14).Value = 1.446177
$ws.Cells.Item(7, 15).Value = 0.7317309745402955
$ws.Cells.Item(7, 16).Value = 0.7317309745402955
$ws.Cells.Item(7, 17).Value = 6.454072470627
$ws.Cells.Item(7, 18).Value = 38.724434823762
$ws.Cells.Item(7, 19).Value = 0.07800407734994777
$ws.Cells.Item(7, 20).Value = 0.07804538973669414

$ws.Cells.Item(8, 1).Value = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value = "Ccl12"
$ws.Cells.Item(8, 3).Value = "Ccr10"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.96682866666667
$ws.Cells.Item(8, 8).Value = 95.900486
$ws.Cells.Item(8, 9).Value = 0.3817886692584505
$ws.Cells.Item(8, 10).Value = 0.3819908715239683
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.5
$ws.Cells.Item(8, 13).Value = 0.2651005
$ws.Cells.Item(8, 14).Value = 0.530201
$ws.Cells.Item(8, 15).Value = 0.2682690254597046
$ws.Cells.Item(8, 16).Value = 0.2682690254597046
$ws.Cells.Item(8, 17).Value = 8.474422262947668
$ws.Cells.Item(8, 18).Value = 50.846533577686
$ws.Cells.Item(8, 19).Value = 0.102422074233522
$ws.Cells.Item(8, 20).Value = 0.1024763188382382

$ws.Cells.Item(9, 1).Value = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value = "Ccl12"
$ws.Cells.Item(9, 3).Value = "Ccr10"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.96682866666667
$ws.Cells.Item(9, 8).Value = 95.900486
$ws.Cells.Item(9, 9).Value = 0.3817886692584505
$ws.Cells.Item(9, 10).Value = 0.3819908715239683
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.5
$ws.Cells.Item(9, 13).Value = 0.7230885
$ws.Cells.Item(9, 14).Value = 1.446177
$ws.Cells.Item(9, 15).Value = 0.7317309745402955
$ws.Cells.Item(9, 16).Value = 0.7317309745402955
$ws.Cells.Item(9, 17).Value = 23.114846190337
$ws.Cells.Item(9, 18).Value = 138.689077142022
$ws.Cells.Item(9, 19).Value = 0.2793665950249286
$ws.Cells.Item(9, 20).Value = 0.2795145526857301

